# Add two new variables ("Body_fat_percentage" and "Exercise_type") as
# extra columns on the "Data" sheet, matching the author's commit:
# "Update example dataset with additional variables".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D1").Value = "Body_fat_percentage"
$ws.Range("E1").Value = "Exercise_type"

# Leave the new columns' data rows empty (only headers were added) and
# select the newly added header cells, matching the saved selection.
$ws.Range("D1:E1").Select() | Out-Null
